$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force-text cells (values that would otherwise be auto-parsed as numbers by Excel)
# are given an explicit "@" (Text) number format before assignment so they are stored
# the same way the source XLSX stores them: as literal text strings.

$ws.Range("D2").Value = '66.030.99'
$ws.Range("E2").Value = '  -3.15%  '
$ws.Range("D3").Value = '3.307.83'
$ws.Range("E3").Value = '  -0.38%  '
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '573.88'
$ws.Range("E5").Value = '  -2.01%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '180.79'
$ws.Range("E6").Value = '  -3.44%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.617'
$ws.Range("E7").Value = '  +2.82%  '
$ws.Range("E8").Value = '  -0.07%  '
$ws.Range("E9").Value = '  -2.48%  '
$ws.Range("E10").Value = '  -0.82%  '
$ws.Range("E11").Value = '  -2.22%  '
$ws.Range("D12").Value = '3.889.28'
$ws.Range("E12").Value = '  -0.31%  '
$ws.Range("E13").Value = '  -1.63%  '
$ws.Range("E14").Value = '  -4.01%  '
$ws.Range("D15").Value = '66.150.06'
$ws.Range("E15").Value = '  -3.34%  '
$ws.Range("E16").Value = '  -1.52%  '
$ws.Range("D17").Value = '3.310.68'
$ws.Range("E17").Value = '  -0.66%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '437.11'
$ws.Range("E18").Value = '  -1.53%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '5.65'
$ws.Range("E19").Value = '  -1.84%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.51'
$ws.Range("E20").Value = '  -1.05%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.54'
$ws.Range("E21").Value = '  -2.79%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '73.25'
$ws.Range("E22").Value = '  -2.71%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.998'
$ws.Range("E23").Value = '  +0.14%  '
$ws.Range("E24").Value = '  -0.11%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0000116'
$ws.Range("E25").Value = '  -3.32%  '
$ws.Range("E26").Value = '  +1.69%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.04'
$ws.Range("E27").Value = '  -2.57%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.999'
$ws.Range("E28").Value = '  +0.54%  '
$ws.Range("E29").Value = '  -2.47%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '22.63'
$ws.Range("E30").Value = '  -1.88%  '
$ws.Range("E31").Value = '  +0.05%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.21'
$ws.Range("E32").Value = '  -3.69%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.73'
$ws.Range("E33").Value = '  -1.99%  '
$ws.Range("E34").Value = '  -4.13%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '159.78'
$ws.Range("E35").Value = '  -2.33%  '
$ws.Range("E36").Value = '  -3.86%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '27.57'
$ws.Range("E37").Value = '  +2.36%  '
$ws.Range("E38").Value = '  -6.02%  '
$ws.Range("D39").Value = '2.827.55'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.787'
$ws.Range("E40").Value = '  -0.40%  '
$ws.Range("E41").Value = '  -3.36%  '
$ws.Range("E42").Value = '  -4.70%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '40.19'
$ws.Range("E43").Value = '  -1.61%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0663'
$ws.Range("E44").Value = '  -2.40%  '
$ws.Range("E45").Value = '  -2.94%  '
$ws.Range("E46").Value = '  -4.56%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '321.75'
$ws.Range("E47").Value = '  -1.83%  '
$ws.Range("E48").Value = '  -2.80%  '
$ws.Range("E49").Value = '  +1.42%  '
$ws.Range("E50").Value = '  -2.10%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.969'
$ws.Range("E51").Value = '  -2.15%  '
